$wb = $excel.ActiveWorkbook

# --- Column width tweaks on "toto" and "tata" sheets ---
$wsToto = $wb.Worksheets.Item("toto")
$wsToto.StandardWidth = 8.23469387755102

$wsTata = $wb.Worksheets.Item("tata")
$wsTata.StandardWidth = 8.23469387755102

# --- "details" sheet: add new value/budget/hours rows ---
$wsDetails = $wb.Worksheets.Item("details")
$wsDetails.Columns.Item(1).ColumnWidth = 8.36734693877551

$wsDetails.Range("D15").Value = "value"
$wsDetails.Range("E15").Value = 50000

$wsDetails.Range("D16").Value = "budget"
$wsDetails.Range("E16").Value = 300

$wsDetails.Range("D19").Value = "hours"
$wsDetails.Range("E19").Value = 85

$wsDetails.Range("E19").Select()
